# Deploy the implementation guide.
# Updates the generation "Date" metadata and adds title-cased "Display"
# values for several concepts in the Concepts sheet, distinct from the
# existing (unchanged) "Code" values.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# 1. Bump the generation timestamp on the Metadata sheet.
$metadata.Range("B8").Value = "2025-09-12T13:34:32+00:00"

# 2. Give each concept its own human-friendly Display text (column C),
#    keeping the machine Code (column B) as-is. Rows 2-5 get a newly
#    title-cased Display string; rows 6-7 already matched, so Display
#    stays equal to Code (but now shares the same string id as Code).
$concepts.Range("C2").Value = "Socially Assigned"
$concepts.Range("C3").Value = "Missing - Restricted Access"
$concepts.Range("C4").Value = "Missing - Not Provided"
$concepts.Range("C5").Value = "Missing - Not Collected"
$concepts.Range("C6").Value = "Derived"
$concepts.Range("C7").Value = "Self-identified"
